# Weekly update: two new price records for Berenjena (Región de Arica y
# Parinacota, fecha 44511) are inserted at the top of the data block,
# pushing the existing rows (formerly 168-194) down to 170-196.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 168.
$ws.Rows("168:169").Insert()

# Row 168: Primera quality record for fecha 44511.
$ws.Range("A168").Value = 9
$ws.Range("B168").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C168").Value = "Metropolitana"
$ws.Range("D168").Value = 44511
$ws.Range("E168").Value = 13
$ws.Range("F168").Value = 100112001
$ws.Range("G168").Value = "Berenjena"
$ws.Range("H168").Value = "Sin especificar"
$ws.Range("I168").Value = "Primera"
$ws.Range("J168").Value = 61
$ws.Range("K168").Value = 9000
$ws.Range("L168").Value = 10000
$ws.Range("M168").Value = 9508
$ws.Range("N168").Value = "`$/caja 60 unidades"
$ws.Range("O168").Value = "Región de Arica y Parinacota"
$ws.Range("P168").Value = 158
$ws.Range("Q168").Value = 60
$ws.Range("R168").Value = "Hortaliza"

# Row 169: Segunda quality record for fecha 44511.
$ws.Range("A169").Value = 9
$ws.Range("B169").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C169").Value = "Metropolitana"
$ws.Range("D169").Value = 44511
$ws.Range("E169").Value = 13
$ws.Range("F169").Value = 100112001
$ws.Range("G169").Value = "Berenjena"
$ws.Range("H169").Value = "Sin especificar"
$ws.Range("I169").Value = "Segunda"
$ws.Range("J169").Value = 25
$ws.Range("K169").Value = 7000
$ws.Range("L169").Value = 7000
$ws.Range("M169").Value = 7000
$ws.Range("N169").Value = "`$/caja 100 unidades"
$ws.Range("O169").Value = "Región de Arica y Parinacota"
$ws.Range("P169").Value = 70
$ws.Range("Q169").Value = 100
$ws.Range("R169").Value = "Hortaliza"
